# Auto-generated edit script: applies value (and where needed, fill/font style)
# corrections to the fund fact-sheet tables across four worksheets, matching the
# updated source data referenced by the commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Dynamic Asset Allocation or Balanced Advantage")
# Row 23
$ws.Range("K23").Value = 68.53
$ws.Range("L23").Value = 31.47
$ws.Range("M23").Value = 18.18

$ws = $wb.Worksheets.Item("Aggressive Hybrid Fund")
# Row 20
$ws.Range("K20").Value = 71.16
$ws.Range("L20").Value = 28.84
$ws.Range("M20").Interior.Color = 8951296
$ws.Range("M20").Font.Color = 0
$ws.Range("M20").Value = 24.72
$ws.Range("N20").Value = 10.83
$ws.Range("P20").Value = 3.23
$ws.Range("Q20").Value = 6.06
$ws.Range("R20").Interior.Color = 14409650
$ws.Range("R20").Font.Color = 0
$ws.Range("R20").Value = 6.49
$ws.Range("S20").Value = 4.12
$ws.Range("T20").Value = 7.43
$ws.Range("U20").Value = 5.24
$ws.Range("V20").Interior.Color = 12959408
$ws.Range("V20").Font.Color = 0
$ws.Range("V20").Value = "-"
$ws.Range("W20").Value = 2.37
$ws.Range("X20").Interior.Color = 15856352
$ws.Range("X20").Font.Color = 0
$ws.Range("X20").Value = 0.67
$ws.Range("Y20").Interior.Color = 11318861
$ws.Range("Y20").Font.Color = 0
$ws.Range("Y20").Value = 17.45
$ws.Range("AB20").Interior.Color = 14409650
$ws.Range("AB20").Font.Color = 0
$ws.Range("AB20").Value = 4.4
# Row 26
$ws.Range("K26").Value = 72.9
$ws.Range("L26").Value = 27.1
$ws.Range("M26").Interior.Color = 8951296
$ws.Range("M26").Font.Color = 0
$ws.Range("M26").Value = 24.44
$ws.Range("N26").Interior.Color = 12897152
$ws.Range("N26").Font.Color = 0
$ws.Range("N26").Value = 10.74
$ws.Range("O26").Value = 1.66
$ws.Range("P26").Value = 1.51
$ws.Range("Q26").Value = 4.62
$ws.Range("R26").Value = 3.96
$ws.Range("S26").Value = 6.52
$ws.Range("T26").Value = 9.29
$ws.Range("U26").Value = 5.61
$ws.Range("W26").Value = 4.54
$ws.Range("Y26").Interior.Color = 12959408
$ws.Range("Y26").Font.Color = 0
$ws.Range("Y26").Value = "-"
$ws.Range("Z26").Value = 15.09
$ws.Range("AB26").Value = 4.52
# Row 27
$ws.Range("K27").Value = 72.9
$ws.Range("L27").Value = 27.1
$ws.Range("M27").Interior.Color = 8951296
$ws.Range("M27").Font.Color = 0
$ws.Range("M27").Value = 24.44
$ws.Range("N27").Interior.Color = 12897152
$ws.Range("N27").Font.Color = 0
$ws.Range("N27").Value = 10.74
$ws.Range("O27").Value = 1.66
$ws.Range("P27").Value = 1.51
$ws.Range("Q27").Value = 4.62
$ws.Range("R27").Value = 3.96
$ws.Range("S27").Value = 6.52
$ws.Range("T27").Value = 9.29
$ws.Range("U27").Value = 5.61
$ws.Range("W27").Value = 4.54
$ws.Range("Y27").Interior.Color = 12959408
$ws.Range("Y27").Font.Color = 0
$ws.Range("Y27").Value = "-"
$ws.Range("Z27").Value = 15.09
$ws.Range("AB27").Value = 4.52
# Row 28
$ws.Range("K28").Value = 72.9
$ws.Range("L28").Value = 27.1
$ws.Range("M28").Interior.Color = 8951296
$ws.Range("M28").Font.Color = 0
$ws.Range("M28").Value = 24.44
$ws.Range("N28").Interior.Color = 12897152
$ws.Range("N28").Font.Color = 0
$ws.Range("N28").Value = 10.74
$ws.Range("O28").Value = 1.66
$ws.Range("P28").Value = 1.51
$ws.Range("Q28").Value = 4.62
$ws.Range("R28").Value = 3.96
$ws.Range("S28").Value = 6.52
$ws.Range("T28").Value = 9.29
$ws.Range("U28").Value = 5.61
$ws.Range("W28").Value = 4.54
$ws.Range("Y28").Interior.Color = 12959408
$ws.Range("Y28").Font.Color = 0
$ws.Range("Y28").Value = "-"
$ws.Range("Z28").Value = 15.09
$ws.Range("AB28").Value = 4.52

$ws = $wb.Worksheets.Item("Conservative Hybrid Fund")
# Row 12
$ws.Range("K12").Value = 24.17
$ws.Range("L12").Value = 75.83
$ws.Range("M12").Value = 8.09
$ws.Range("N12").Value = 3.11
$ws.Range("P12").Value = 0.95
$ws.Range("Q12").Value = 1.98
$ws.Range("R12").Value = 2.26
$ws.Range("S12").Value = 1.74
$ws.Range("T12").Value = 1.81
$ws.Range("U12").Value = 1.89
$ws.Range("V12").Interior.Color = 12959408
$ws.Range("V12").Font.Color = 0
$ws.Range("V12").Value = "-"
$ws.Range("W12").Value = 1.87
$ws.Range("X12").Value = 0.45
$ws.Range("Y12").Interior.Color = 12897152
$ws.Range("Y12").Font.Color = 0
$ws.Range("Y12").Value = 9.61
$ws.Range("AB12").Value = 47.62
$ws.Range("AC12").Interior.Color = 15856352
$ws.Range("AC12").Font.Color = 0
$ws.Range("AC12").Value = 2.85
$ws.Range("AE12").Interior.Color = 14409650
$ws.Range("AE12").Font.Color = 0
$ws.Range("AE12").Value = 3.66
# Row 15
$ws.Range("K15").Interior.Color = 10135078
$ws.Range("K15").Font.Color = 0
$ws.Range("K15").Value = 21.63
$ws.Range("L15").Value = 78.37
$ws.Range("M15").Interior.Color = 14409650
$ws.Range("M15").Font.Color = 0
$ws.Range("M15").Value = 6.9
$ws.Range("N15").Value = 2.25
$ws.Range("O15").Value = 1.25
$ws.Range("P15").Value = 0.55
$ws.Range("Q15").Value = 1.46
$ws.Range("R15").Value = 0.98
$ws.Range("S15").Value = 1.62
$ws.Range("T15").Value = 3.87
$ws.Range("U15").Value = 1.59
$ws.Range("V15").Interior.Color = 8023636
$ws.Range("V15").Font.Color = 0
$ws.Range("V15").Value = 0
$ws.Range("W15").Value = 1.17
$ws.Range("Y15").Interior.Color = 8951296
$ws.Range("Y15").Font.Color = 0
$ws.Range("Y15").Value = 25.52
$ws.Range("AB15").Interior.Color = 6056192
$ws.Range("AB15").Font.Color = 15855596
$ws.Range("AB15").Value = 42.67

$ws = $wb.Worksheets.Item("Multi Asset Allocation")
# Row 10
$ws.Range("G10").Value = 18.57
# Row 12
$ws.Range("G12").Value = 10.59

